$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.057.72"
$ws.Range("E2").Value = "  +1.27%  "

$ws.Range("D3").Value = "'3.777.53"
$ws.Range("E3").Value = "  -0.67%  "

$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").Value = "'626.79"
$ws.Range("E5").Value = "  +4.22%  "

$ws.Range("D6").Value = "'163.79"
$ws.Range("E6").Value = "  -1.02%  "

$ws.Range("D7").Value = "'3.772.33"
$ws.Range("E7").Value = "  -0.72%  "

$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("D9").Value = "'0.519"
$ws.Range("E9").Value = "  -0.03%  "

$ws.Range("E10").Value = "  +0.54%  "

$ws.Range("D11").Value = "'0.450"
$ws.Range("E11").Value = "  -0.43%  "

$ws.Range("D12").Value = "'6.62"
$ws.Range("E12").Value = "  +2.24%  "

$ws.Range("E13").Value = "  -1.31%  "

$ws.Range("D14").Value = "'35.43"
$ws.Range("E14").Value = "  -1.37%  "

$ws.Range("D15").Value = "'4.419.79"
$ws.Range("E15").Value = "  -0.41%  "

$ws.Range("D16").Value = "'3.816.72"
$ws.Range("E16").Value = "  +0.27%  "

$ws.Range("D17").Value = "'69.038.44"
$ws.Range("E17").Value = "  +1.28%  "

$ws.Range("D18").Value = "'17.92"
$ws.Range("E18").Value = "  -2.96%  "

$ws.Range("E19").Value = "  -1.04%  "

$ws.Range("D20").Value = "'7.07"
$ws.Range("E20").Value = "  -0.58%  "

$ws.Range("D21").Value = "'467.10"
$ws.Range("E21").Value = "  +1.04%  "

$ws.Range("D22").Value = "'9.62"
$ws.Range("E22").Value = "  -1.19%  "

$ws.Range("E23").Value = "  +0.11%  "

$ws.Range("D24").Value = "'0.0000149"
$ws.Range("E24").Value = "  -0.56%  "

$ws.Range("D25").Value = "'83.07"
$ws.Range("E25").Value = "  -0.01%  "

$ws.Range("D26").Value = "'11.97"
$ws.Range("E26").Value = "  -0.90%  "

$ws.Range("E27").Value = "  +0.96%  "

$ws.Range("E28").Value = "  +0.00%  "

$ws.Range("E29").Value = "  -0.35%  "

$ws.Range("D30").Value = "'3.933.23"
$ws.Range("E30").Value = "  -0.41%  "

$ws.Range("E31").Value = "  +0.32%  "

$ws.Range("E32").Value = "  -1.28%  "

$ws.Range("D33").Value = "'7.23"
$ws.Range("E33").Value = "  -1.65%  "

$ws.Range("D34").Value = "'28.88"
$ws.Range("E34").Value = "  -1.94%  "

$ws.Range("E35").Value = "  -0.04%  "

$ws.Range("D36").Value = "'3.721.69"
$ws.Range("E36").Value = "  -0.73%  "

$ws.Range("E37").Value = "  -1.03%  "

$ws.Range("E38").Value = "  +2.40%  "

$ws.Range("E39").Value = "  +7.66%  "

$ws.Range("E40").Value = "  -0.86%  "

$ws.Range("E41").Value = "  +0.26%  "

$ws.Range("D42").Value = "'0.971"
$ws.Range("E42").Value = "  -1.92%  "

$ws.Range("E43").Value = "  +0.25%  "

$ws.Range("E44").Value = "  +0.01%  "

$ws.Range("D45").Value = "'153.59"
$ws.Range("E45").Value = "  +1.28%  "

$ws.Range("E46").Value = "  -1.38%  "

$ws.Range("D47").Value = "'46.86"
$ws.Range("E47").Value = "  -1.25%  "

$ws.Range("D48").Value = "'1.91"
$ws.Range("E48").Value = "  +1.85%  "

$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "'8.40"
$ws.Range("E49").Value = "  +0.24%  "

$ws.Range("B50").Value = "Arweave"
$ws.Range("C50").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D50").Value = "'42.20"
$ws.Range("E50").Value = "  -2.14%  "

$ws.Range("E51").Value = "  +1.77%  "
